# nuevas funcionalidades y actualización
# Add the new match result row (row 14) to the Spartak Varna stats sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "24/10/2025"
$ws.Range("B14").Value = "Spartak Varna"
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "Botev Plovdiv"
$ws.Range("F14").Value = "W"
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 1
$ws.Range("M14").Value = 11
$ws.Range("N14").Value = 11
$ws.Range("O14").Value = 7
$ws.Range("P14").Value = 5
